# Supporting_Information.xlsx - "Add files via upload" edit
#
# 1. Rename the four data-tab worksheets, dropping the "I" from the
#    "SIx_" prefix (SI1 -> S1, SI2 -> S2, SI3 -> S3) and renumbering the
#    last one from "SI3_Riverton_Obs" to "S4_Riverton_Obs".
# 2. Update the "Tab Name" column on the ReadMe sheet to match (note the
#    last entry becomes "I4_Riverton_Obs" - a different typo than the
#    sheet tab got).
# 3. Move the active sheet/selection from S4_Riverton_Obs (H17) to
#    S3_Synthetic_Test_Obs (H31), and update the ReadMe selection too.

$wb = $excel.ActiveWorkbook

$wsReadMe = $wb.Worksheets.Item(1)
$wsS1 = $wb.Worksheets.Item(2)
$wsS2 = $wb.Worksheets.Item(3)
$wsS3 = $wb.Worksheets.Item(4)
$wsS4 = $wb.Worksheets.Item(5)

# Rename the worksheet tabs.
$wsS1.Name = "S1_Synthetic_Test_Tracer_Input"
$wsS2.Name = "S2_Riverton_Tracer_Input"
$wsS3.Name = "S3_Synthetic_Test_Obs"
$wsS4.Name = "S4_Riverton_Obs"

# Update the matching labels in the ReadMe table.
$wsReadMe.Range("A2").Value = "S1_Synthetic_Test_Tracer_Input"
$wsReadMe.Range("A3").Value = "S2_Riverton_Tracer_Input"
$wsReadMe.Range("A4").Value = "S3_Synthetic_Test_Obs"
$wsReadMe.Range("A5").Value = "I4_Riverton_Obs"

# Update the ReadMe sheet's own selection (not the active sheet in the end).
$wsReadMe.Range("A5").Select()

# Make S3_Synthetic_Test_Obs the active sheet/selection, S4_Riverton_Obs
# loses tabSelected and keeps its own prior selection.
$wsS3.Activate()
$wsS3.Range("H31").Select()
